$d = $word.ActiveDocument

$ndash = [char]0x2013
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Remove the w:proofErr (spellStart/spellEnd) markers that wrap "DebugUni"
#    in the very first paragraph. The cleanest way to drop those proofing
#    markers is to replace the whole paragraph (pPr + runs) with a fresh
#    copy of itself built from raw WordprocessingML, which never contains
#    the proofErr markers.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstParaRange = $d.Range($firstPara.Range.Start, $firstPara.Range.End)
[void]$firstParaRange.Delete()

$debugUniXml = "<w:p $wns>" +
    "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr>" +
    "<w:r><w:t>DebugUni</w:t></w:r>" +
    "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
    "</w:p>"

$insertionPoint = $d.Range(0, 0)
[void]$insertionPoint.InsertXML($debugUniXml)

# ---------------------------------------------------------------------------
# 2) Add two new suggestions (three new list paragraphs in total) after the
#    "Работа за пичове" bullet. The hidden _GoBack bookmark needs to end up
#    wrapping the very end of the new last paragraph, so we delete it from
#    its old location and re-create it inside the XML of the new last
#    paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
[void]$goBack.Delete()

$newParasXml =
    "<w:p $wns>" +
        "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr>" +
        "<w:r><w:t>WizardUni $ndash want to be a wizard of programming? Apply!</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wns>" +
        "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr></w:pPr>" +
        "<w:r><w:t>Use wizard image design from Softuni presentation materials.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wns>" +
        "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr></w:pPr>" +
        "<w:r><w:t>Site background $ndash use the binary code design in Softuni classrooms.</w:t></w:r>" +
        "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>" +
    "</w:p>"

$endOfDoc = $d.Range($d.Content.End, $d.Content.End)
[void]$endOfDoc.InsertXML($newParasXml)
